$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Total" header column (T1) ---
$ws.Range("T1").Value = "Total"

# --- Existing category rows (2-6) keep their text; just add the new Total column (T) ---
$ws.Range("T2").Value = 79933
$ws.Range("T3").Value = 6878
$ws.Range("T4").Value = 32858
$ws.Range("T5").Value = 13198
$ws.Range("T6").Value = 46288

# --- New row 7: "Outros" category ---
$ws.Range("A7").Value = "Outros"
$ws.Range("B7").Value = 7343
$ws.Range("C7").Value = 409
$ws.Range("D7").Value = 655
$ws.Range("E7").Value = 2147
$ws.Range("F7").Value = 3526
$ws.Range("G7").Value = 3700
$ws.Range("H7").Value = 3922
$ws.Range("I7").Value = 4200
$ws.Range("J7").Value = 4694
$ws.Range("K7").Value = 5177
$ws.Range("L7").Value = 5565
$ws.Range("M7").Value = 5367
$ws.Range("N7").Value = 4976
$ws.Range("O7").Value = 4882
$ws.Range("P7").Value = 5506
$ws.Range("Q7").Value = 6043
$ws.Range("R7").Value = 17027
$ws.Range("S7").Value = 657
$ws.Range("T7").Value = 85796

# --- New row 8: "Total" row (grand total across all categories) ---
$ws.Range("A8").Value = "Total"
$ws.Range("B8").Value = 8274
$ws.Range("C8").Value = 625
$ws.Range("D8").Value = 920
$ws.Range("E8").Value = 2632
$ws.Range("F8").Value = 4207
$ws.Range("G8").Value = 4818
$ws.Range("H8").Value = 5491
$ws.Range("I8").Value = 6702
$ws.Range("J8").Value = 8902
$ws.Range("K8").Value = 12427
$ws.Range("L8").Value = 16163
$ws.Range("M8").Value = 19139
$ws.Range("N8").Value = 21391
$ws.Range("O8").Value = 23257
$ws.Range("P8").Value = 27274
$ws.Range("Q8").Value = 29766
$ws.Range("R8").Value = 72121
$ws.Range("S8").Value = 842
$ws.Range("T8").Value = 264951
